# Simulation and efficiency update: dedupe player rows, recompute pl21 minutes (reset to 0
# for players with no 2021 minutes) and re-derive the running max(pl20, pl21) column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old trailing rows (27-31) -- the rebuilt table only needs rows 1-26.
$ws.Range("A27:D31").EntireRow.Delete()

$rows = @(
    @("player", "pl20_minutes", "pl21_minutes", "max"),
    @("Daniel James", 910, 0, 910),
    @("Dean Henderson", 1125, 0, 1125),
    @("Eric Bailly", 916, 0, 916),
    @("Anthony Martial", 1480, 0, 1480),
    @("Donny van de Beek", 524, 0, 524),
    @("Juan Mata", 509, 0, 509),
    @("David de Gea", 2295, 3420, 3420),
    @("Bruno Fernandes", 3099, 3111, 3111),
    @("Harry Maguire", 3047, 2515, 3047),
    @("Cristiano Ronaldo", 0, 2456, 2456),
    @("Scott McTominay", 2129, 2387, 2387),
    @("Victor LindelÃ¶f", 2585, 2355, 2585),
    @("Fred", 2390, 2028, 2390),
    @("Jadon Sancho", 0, 1900, 1900),
    @("RaphaÃ«l Varane", 0, 1831, 1831),
    @("Diogo Dalot", 0, 1801, 1801),
    @("Aaron Wan-Bissaka", 3060, 1793, 3060),
    @("Alex Telles", 690, 1660, 1660),
    @("Luke Shaw", 2654, 1601, 2654),
    @("Nemanja MatiÄ‡", 1106, 1382, 1382),
    @("Paul Pogba", 1897, 1349, 1897),
    @("Mason Greenwood", 1822, 1271, 1822),
    @("Marcus Rashford", 2920, 1235, 2920),
    @("Anthony Elanga", 0, 1214, 1214),
    @("Edinson Cavani", 1375, 766, 1375),
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 1, $j + 1).Value = $row[$j]
    }
}

$ws.Range("G7").Select()

